# Extend the "2024-1" tracker sheet with one more week of rows (53-59),
# mirroring the pattern already used for the previous weeks:
#   - "empty" placeholder days (only ID + DATE filled in) get the same
#     formatting as the existing placeholder rows (e.g. row 48)
#   - the one day that has an actual bet (row 59) gets filled in fully,
#     following the same formulas as the previous data rows (e.g. row 52)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024-1")

# ---------------------------------------------------------------------
# 1) Six "no bet placed" days: 2024-02-05 .. 2024-02-10 (rows 53-58)
# ---------------------------------------------------------------------
$newRows   = @(53, 54, 55, 56, 57, 58)
$idValues  = @(52, 53, 54, 55, 56, 57)
$dateSerials = @(45327, 45328, 45329, 45330, 45331, 45332)

for ($k = 0; $k -lt $newRows.Length; $k++) {
    $r = $newRows[$k]

    # Values first (so the subsequent format paste below doesn't get
    # overridden by any auto number-format inheritance)
    $ws.Range("A$r").Value = $idValues[$k]
    $ws.Range("C$r").Value = $dateSerials[$k]

    # Copy the look of the previous empty-day row (48) onto the new one,
    # cell-by-cell so we don't stamp unrelated/blank columns with styles
    # they never had.
    $ws.Range("A48").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    $ws.Range("C48").Copy() | Out-Null
    $ws.Range("C$r").PasteSpecial(-4122) | Out-Null

    $ws.Range("H48:K48").Copy() | Out-Null
    $ws.Range("H$r`:K$r").PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------
# 2) The new bet on 2024-02-11 (row 59)
# ---------------------------------------------------------------------
$ws.Range("A59").Value = 58
$ws.Range("B59").Value = 39
$ws.Range("C59").Value = 45333
$ws.Range("D59").Value = 1.313
$ws.Range("E59").Value = 1
$ws.Range("F59").Value = 2584

$ws.Range("G59").Formula = "=F59*E59*D59"
$ws.Range("H59").Formula = "=I52"
$ws.Range("I59").Formula = "=H59+G59-F59"
$ws.Range("J59").Formula = "=I59-H59"
$ws.Range("K59").Formula = "=I59/`$H`$2-1"

$ws.Range("L59").Value = "EUROPE"
$ws.Range("M59").Value = "ESPORTS"
$ws.Range("N59").Value = "LEC"

# Mirror the formatting of the previous full data row (52) onto row 59.
$ws.Range("A52:N52").Copy() | Out-Null
$ws.Range("A59:N59").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3) Scroll/selection bookkeeping, matching where the user ended up
#    after typing in the new rows.
# ---------------------------------------------------------------------
$ws.Activate()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 43
$aw.ScrollColumn = 2
$ws.Range("O59").Select()
